$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = "Risoluzione bug critico su fetch prezzi yfinance. 1) I titoli USA restituivano NaN per la giornata corrente (pre-market/festivi). 2) La logica batch aveva un duplicato e gestiva male il MultiIndex per singoli ticker."
$ws.Range("E7").Value = "Analisi log, creazione script di riproduzione (verify_nan_issue.py, debug_backend_logic.py), patch logica backend."
$ws.Range("F7").Value = "backend/main.py"
$ws.Range("I7").Value = "NO"
$ws.Range("J7").Value = "Data Consistency, Real-time Pricing"
$ws.Range("K7").Value = "SI"
$ws.Range("L7").Value = "Verifica caricamento Portfolio con mix titoli EU/USA; Verifica fallback."
$ws.Range("M7").Value = "Prezzi corretti visualizzati in dashboard. Fallback attivo."

$wb.Save()
